$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93 (new data point for 2021-02-08), pushing the
# previously existing rows 93-113 down to 94-114.
$ws.Rows(93).Insert()

# The insert leaves the new row's A cell with a generic style; restore the
# date-column formatting by copying it down from the row above.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)

# New row 93 data: 2021-02-08 (serial 44235), 19 new cases that day.
$ws.Cells.Item(93, 1).Value = 44235
$ws.Cells.Item(93, 2).Value = 19
$ws.Cells.Item(93, 3).Value = 60
$ws.Cells.Item(93, 4).Value = 148.964695367198

# Rolling 7-day totals (col C) / per-100k rates (col D) change for the rows
# around the inserted point because the trailing window now includes it.
$ws.Cells.Item(90, 3).Value = 62
$ws.Cells.Item(90, 4).Value = 153.9301852127712
$ws.Cells.Item(91, 3).Value = 65
$ws.Cells.Item(91, 4).Value = 161.3784199811311
$ws.Cells.Item(92, 3).Value = 66
$ws.Cells.Item(92, 4).Value = 163.8611649039178

$ws.Cells.Item(94, 3).Value = 68
$ws.Cells.Item(94, 4).Value = 168.826654749491
$ws.Cells.Item(95, 3).Value = 73
$ws.Cells.Item(95, 4).Value = 181.2403793634242
$ws.Cells.Item(96, 3).Value = 66
$ws.Cells.Item(96, 4).Value = 163.8611649039178

$ws.Cells.Item(111, 3).Value = 224
$ws.Cells.Item(111, 4).Value = 556.1348627042058
$ws.Cells.Item(112, 3).Value = 224
$ws.Cells.Item(112, 4).Value = 556.1348627042058

# Append a brand-new final row (2021-03-02, serial 44257, 25 new cases) after
# the last existing (now-shifted) row.
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 25

# Columns C/D on the trailing rows (113-115) stay blank, same as the other
# not-yet-computed rows -- write them as empty text (matching the existing
# blank cells) rather than leaving them as untyped/empty numeric cells.
$ws.Cells.Item(115, 3).Value = "'"
$ws.Cells.Item(115, 4).Value = "'"
$ws.Cells.Item(115, 3).ClearFormats()
$ws.Cells.Item(115, 4).ClearFormats()
